$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3-10 down to 4-11
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with data
$ws.Cells.Item(3, 1).Value = 7
$ws.Cells.Item(3, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(3, 3).Value = "Ñuble"
$ws.Cells.Item(3, 4).Value = 44453
$ws.Cells.Item(3, 4).NumberFormat = $ws.Cells.Item(4, 4).NumberFormat
$ws.Cells.Item(3, 5).Value = 16
$ws.Cells.Item(3, 6).Value = 100112013
$ws.Cells.Item(3, 7).Value = "Alcachofa"
$ws.Cells.Item(3, 8).Value = "Madrigal"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 160
$ws.Cells.Item(3, 11).Value = 12500
$ws.Cells.Item(3, 12).Value = 13000
$ws.Cells.Item(3, 13).Value = 12750
$ws.Cells.Item(3, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(3, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(3, 16).Value = 319
$ws.Cells.Item(3, 17).Value = 40
$ws.Cells.Item(3, 18).Value = "Hortaliza"
